$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 38776.652
$ws.Range("I53").Value = 71775.71000000001
$ws.Range("J53").Value = 277.75
$ws.Range("K53").Value = 71775.71000000001
$ws.Range("L53").Value = 277.75
$ws.Range("M53").Value = -71138.71000000001
$ws.Range("N53").Value = -1551.75
$ws.Range("H94").Value = 5534.3335
$ws.Range("I94").Value = 5534.3335
$ws.Range("K94").Value = 5534.3335
$ws.Range("M94").Value = -5083.3335
$ws.Range("H125").Value = 2279.7058
$ws.Range("J125").Value = 2156.6365
$ws.Range("L125").Value = 19409.7285
$ws.Range("N125").Value = -24329.7285
$ws.Range("H129").Value = 1502.4
$ws.Range("I129").Value = 720.75
$ws.Range("J129").Value = 2023.5
$ws.Range("K129").Value = 2162.25
$ws.Range("L129").Value = 6070.5
$ws.Range("M129").Value = 2837.75
$ws.Range("N129").Value = -16070.5
$ws.Range("H137").Value = 3743
$ws.Range("I137").Value = 2725.3225
$ws.Range("J137").Value = 6611
$ws.Range("K137").Value = 8175.967500000001
$ws.Range("L137").Value = 19833
$ws.Range("M137").Value = -5625.967500000001
$ws.Range("N137").Value = -24933
$ws.Range("H138").Value = 5566.615
$ws.Range("I138").Value = 1970.1428
$ws.Range("J138").Value = 6891.6313
$ws.Range("K138").Value = 5910.428400000001
$ws.Range("L138").Value = 20674.8939
$ws.Range("M138").Value = -770.4284000000007
$ws.Range("N138").Value = -30954.8939
$ws.Range("H141").Value = 2541.1667
$ws.Range("I141").Value = 2353
$ws.Range("J141").Value = 2998.1428
$ws.Range("K141").Value = 7059
$ws.Range("L141").Value = 8994.428400000001
$ws.Range("M141").Value = -1879
$ws.Range("N141").Value = -19354.4284

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 25000076
$ws.Range("I10").Value = 25000076
$ws.Range("K10").Value = 25000076
$ws.Range("M10").Value = -24999906
$ws.Range("H32").Value = 712.66
$ws.Range("I32").Value = 617.3262999999999
$ws.Range("J32").Value = 2524
$ws.Range("K32").Value = 617.3262999999999
$ws.Range("L32").Value = 2524
$ws.Range("M32").Value = -330.3262999999999
$ws.Range("N32").Value = -3098
$ws.Range("H132").Value = 27152.334
$ws.Range("I132").Value = 2201.2
$ws.Range("J132").Value = 151908
$ws.Range("K132").Value = 6603.599999999999
$ws.Range("L132").Value = 455724
$ws.Range("M132").Value = -4073.599999999999
$ws.Range("N132").Value = -460784

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 10569.5
$ws.Range("J81").Value = 10569.5
$ws.Range("L81").Value = 10569.5
$ws.Range("N81").Value = -12691.5
$ws.Range("H84").Value = 10569.5
$ws.Range("J84").Value = 10569.5
$ws.Range("L84").Value = 31708.5
$ws.Range("N84").Value = -42316.5
$ws.Range("H140").Value = 77834.414
$ws.Range("J140").Value = 82183
$ws.Range("L140").Value = 82183
$ws.Range("N140").Value = -92543

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 703.1667
$ws.Range("J5").Value = 1198
$ws.Range("L5").Value = 1198
$ws.Range("N5").Value = -1422
$ws.Range("H31").Value = 3962.4849
$ws.Range("I31").Value = 1136.75
$ws.Range("J31").Value = 6622
$ws.Range("K31").Value = 1136.75
$ws.Range("L31").Value = 6622
$ws.Range("M31").Value = -841.75
$ws.Range("N31").Value = -7212
$ws.Range("H34").Value = 3962.4849
$ws.Range("I34").Value = 1136.75
$ws.Range("J34").Value = 6622
$ws.Range("K34").Value = 1136.75
$ws.Range("L34").Value = 6622
$ws.Range("M34").Value = -934.75
$ws.Range("N34").Value = -7026
$ws.Range("H134").Value = 15589.046
$ws.Range("I134").Value = 11599.2
$ws.Range("J134").Value = 16762.53
$ws.Range("K134").Value = 34797.60000000001
$ws.Range("L134").Value = 50287.59
$ws.Range("M134").Value = -32262.60000000001
$ws.Range("N134").Value = -55357.59
$ws.Range("H140").Value = 11400
$ws.Range("J140").Value = 11400
$ws.Range("L140").Value = 11400
$ws.Range("N140").Value = -21760

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9021203
$ws.Range("I4").Value = 12135858
$ws.Range("K4").Value = 36407574
$ws.Range("M4").Value = -36407462
$ws.Range("H5").Value = 114571.34
$ws.Range("I5").Value = 854.4737
$ws.Range("J5").Value = 200996.16
$ws.Range("K5").Value = 2563.4211
$ws.Range("L5").Value = 602988.48
$ws.Range("M5").Value = -2451.4211
$ws.Range("N5").Value = -603212.48
$ws.Range("H37").Value = 141241
$ws.Range("J37").Value = 141241
$ws.Range("L37").Value = 423723
$ws.Range("N37").Value = -423947
$ws.Range("H107").Value = 4287.9688
$ws.Range("I107").Value = 476.83334
$ws.Range("J107").Value = 4682.224
$ws.Range("K107").Value = 1430.50002
$ws.Range("L107").Value = 14046.672
$ws.Range("M107").Value = 489.4999800000001
$ws.Range("N107").Value = -17886.672
$ws.Range("H122").Value = 114075.445
$ws.Range("J122").Value = 134973.58
$ws.Range("L122").Value = 1214762.22
$ws.Range("N122").Value = -1219662.22
$ws.Range("H124").Value = 11858.167
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 11858.167
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 35574.501
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -45394.501
$ws.Range("H135").Value = 114571.34
$ws.Range("I135").Value = 854.4737
$ws.Range("J135").Value = 200996.16
$ws.Range("K135").Value = 7690.263300000001
$ws.Range("L135").Value = 1808965.44
$ws.Range("M135").Value = -5155.263300000001
$ws.Range("N135").Value = -1814035.44

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 6251250
$ws.Range("I3").Value = 8334500
$ws.Range("J3").Value = 1499.5
$ws.Range("K3").Value = 8334500
$ws.Range("L3").Value = 1499.5
$ws.Range("M3").Value = -8334384
$ws.Range("N3").Value = -1731.5
$ws.Range("H11").Value = 6793721
$ws.Range("J11").Value = 2208897.8
$ws.Range("L11").Value = 2208897.8
$ws.Range("N11").Value = -2209175.8
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H126").Value = 5744.8965
$ws.Range("I126").Value = 2584.6843
$ws.Range("J126").Value = 11749.3
$ws.Range("K126").Value = 7754.0529
$ws.Range("L126").Value = 35247.89999999999
$ws.Range("M126").Value = -5284.0529
$ws.Range("N126").Value = -40187.89999999999
$ws.Range("H132").Value = 3454.68
$ws.Range("I132").Value = 2500.8948
$ws.Range("J132").Value = 6475
$ws.Range("K132").Value = 7502.6844
$ws.Range("L132").Value = 19425
$ws.Range("M132").Value = -4972.6844
$ws.Range("N132").Value = -24485

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6908.522
$ws.Range("I7").Value = 4357.6665
$ws.Range("J7").Value = 11691.375
$ws.Range("K7").Value = 4357.6665
$ws.Range("L7").Value = 11691.375
$ws.Range("M7").Value = -4245.6665
$ws.Range("N7").Value = -11915.375
$ws.Range("H40").Value = 4186.84
$ws.Range("I40").Value = 2568.5881
$ws.Range("J40").Value = 7625.625
$ws.Range("K40").Value = 2568.5881
$ws.Range("L40").Value = 7625.625
$ws.Range("M40").Value = -2432.5881
$ws.Range("N40").Value = -7897.625
$ws.Range("H126").Value = 6908.522
$ws.Range("I126").Value = 4357.6665
$ws.Range("J126").Value = 11691.375
$ws.Range("K126").Value = 13072.9995
$ws.Range("L126").Value = 35074.125
$ws.Range("M126").Value = -10602.9995
$ws.Range("N126").Value = -40014.125

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H132").Value = 4875.098
$ws.Range("I132").Value = 4474.814
$ws.Range("J132").Value = 7026.625
$ws.Range("K132").Value = 13424.442
$ws.Range("L132").Value = 21079.875
$ws.Range("M132").Value = -10894.442
$ws.Range("N132").Value = -26139.875
